$wb = $excel.ActiveWorkbook

$weekLabel = "24/01/2022 - 30/01/2022"

# --- Infanzia ---
$ws = $wb.Worksheets.Item("Infanzia")
$ws.Range("A13").Value = $weekLabel
$ws.Range("B13").Value = 29
$ws.Range("C13").Value = 189
$ws.Range("D13").Value = 218
$ws.Range("E13").Select()

# --- Primaria ---
$ws = $wb.Worksheets.Item("Primaria")
$ws.Range("A13").Value = $weekLabel
$ws.Range("B13").Value = 35
$ws.Range("C13").Value = 445
$ws.Range("D13").Value = 480
$ws.Range("E13").Select()

# --- Media ---
$ws = $wb.Worksheets.Item("Media")
$ws.Range("A13").Value = $weekLabel
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 166
$ws.Range("D13").Value = 170
$ws.Range("E13").Select()

# --- Superiore ---
$ws = $wb.Worksheets.Item("Superiore")
$ws.Range("A13").Value = $weekLabel
$ws.Range("B13").Value = 16
$ws.Range("C13").Value = 226
$ws.Range("D13").Value = 242
$ws.Range("E13").Select()

# --- Totale casi (active sheet) ---
$ws = $wb.Worksheets.Item("Totale casi")
$ws.Range("A13").Value = $weekLabel
$ws.Range("B13").Value = 84
$ws.Range("C13").Value = 1026
$ws.Range("D13").Value = 1110
$ws.Activate()
$ws.Range("B14").Select()
